$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 58:59 - this shifts the existing data (old rows
# 58..207) down to 60..209, which also naturally produces the new trailing
# rows 208/209 as duplicates of the former 206/207 (matching the target diff).
$ws.Rows("58:59").Insert()

# Row 58 - new "Primera" quality record for date 44544 (2021-12-14)
$ws.Cells.Item(58, 1).Value = 8
$ws.Cells.Item(58, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(58, 3).Value = "Coquimbo"
$ws.Cells.Item(58, 4).Value = 44544
$ws.Cells.Item(58, 5).Value = 4
$ws.Cells.Item(58, 6).Value = 100114014
$ws.Cells.Item(58, 7).Value = "Betarraga"
$ws.Cells.Item(58, 8).Value = "Sin especificar"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 3100
$ws.Cells.Item(58, 11).Value = 450
$ws.Cells.Item(58, 12).Value = 500
$ws.Cells.Item(58, 13).Value = 475
$ws.Cells.Item(58, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(58, 15).Value = "Provincia del Elqu$([char]0x00ED)"
$ws.Cells.Item(58, 16).Value = 158
$ws.Cells.Item(58, 17).Value = 3
$ws.Cells.Item(58, 18).Value = "Hortaliza"

# Row 59 - new "Segunda" quality record for date 44544 (2021-12-14)
$ws.Cells.Item(59, 1).Value = 8
$ws.Cells.Item(59, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(59, 3).Value = "Coquimbo"
$ws.Cells.Item(59, 4).Value = 44544
$ws.Cells.Item(59, 5).Value = 4
$ws.Cells.Item(59, 6).Value = 100114014
$ws.Cells.Item(59, 7).Value = "Betarraga"
$ws.Cells.Item(59, 8).Value = "Sin especificar"
$ws.Cells.Item(59, 9).Value = "Segunda"
$ws.Cells.Item(59, 10).Value = 1520
$ws.Cells.Item(59, 11).Value = 350
$ws.Cells.Item(59, 12).Value = 400
$ws.Cells.Item(59, 13).Value = 375
$ws.Cells.Item(59, 14).Value = "`$/paquete 3 unidades"
$ws.Cells.Item(59, 15).Value = "Provincia del Elqu$([char]0x00ED)"
$ws.Cells.Item(59, 16).Value = 125
$ws.Cells.Item(59, 17).Value = 3
$ws.Cells.Item(59, 18).Value = "Hortaliza"
